$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 759, pushing the existing rows 759-800 down to 760-801.
$ws.Rows(759).Insert()

# Populate the newly inserted row with the new data point for 2026/02/05.
# Force the date column to text format first so Excel doesn't auto-convert
# the "yyyy/mm/dd" string into a date serial number (matches the other
# rows, which all store the date as literal text).
$ws.Range("A759").NumberFormat = "@"
$ws.Range("A759").Value = "2026/02/05"
$ws.Range("B759").Value = "木"
$ws.Range("C759").Value = 14
$ws.Range("D759").Value = 60
